$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# ---------------------------------------------------------------------------
# The Cypher queries stored in the "startup" sheet were rewritten (the
# participants/samples/files per-tab queries now route through a common
# CALL{...} sub-query pattern with OPTIONAL MATCH / diagnosis joins, and the
# combined "StatQuery" cell was rebuilt around three CALL blocks). The two
# output-file-name columns (Neo4jData.xlsx / WebData.xlsx) keep the same
# text, only the query text actually changes.
# ---------------------------------------------------------------------------

$neo4jFile = "TC04_CDS_Filter_RefGenAssembly-NotAvailable_Neo4jData.xlsx"
$webFile   = "TC04_CDS_Filter_RefGenAssembly-NotAvailable_WebData.xlsx"

$statQuery = @"
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in ['Not Available']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in ['Not Available']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in ['Not Available']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS ``Files``
"@

$participantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in ['Not Available'] 
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id limit 100
"@

$filesQuery = @"
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in ['Not Available']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name,'') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id, '') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
ORDER BY f.file_name limit 100
"@

$samplesQuery = @"
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in ['Not Available']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(p.participant_id,'') as ``Participant ID``,
    coalesce(s.study_name, '') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(samp.sample_tumor_status,'') as ``Tumor``,
    coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER BY samp.sample_id limit 100
"@

# Fill column-by-column (D, then E, then C, then the per-tab B cells) so the
# shared-string table is rebuilt in the same order as the source edit.
$ws.Range("D2").Value = $neo4jFile
$ws.Range("D3").Value = $neo4jFile
$ws.Range("D4").Value = $neo4jFile

$ws.Range("E2").Value = $webFile
$ws.Range("E3").Value = $webFile
$ws.Range("E4").Value = $webFile

$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# Row 2 = ParticipantsTab
$ws.Range("B2").Value = $participantsQuery

# Row 4 = FilesTab
$ws.Range("B4").Value = $filesQuery

# Row 3 = SamplesTab
$ws.Range("B3").Value = $samplesQuery

# Row heights reflow with the new (differently-wrapped) query text.
$ws.Rows.Item(2).RowHeight = 242.25
$ws.Rows.Item(3).RowHeight = 260.25
$ws.Rows.Item(4).RowHeight = 279.75
